# Swap the order of the two recorder names within the "Recorded By"
# column (column G) of the "Session Analysis Results" worksheet.
#
#   Before: "System, dnasr281@gmail.com"
#   After:  "dnasr281@gmail.com, System"
#
# Every row whose column-G value is exactly the old combined string is
# updated; rows holding only "System" or only "dnasr281@gmail.com" (or
# anything else) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
